$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Section_A and Section_B timetables: tutorial slots changed from 1.5hr
# blocks to 1hr blocks, which means the grid now needs more rows to cover
# the same overall day span (12:00-18:30). Extend both timetables with the
# new 1-hour rows (times 12:00-13:00, 13:00-14:00, 15:30-16:30, 16:30-17:30,
# 17:30-18:30), all still "Free".
# ---------------------------------------------------------------------------
$newTimes = @("12:00-13:00", "13:00-14:00", "15:30-16:30", "16:30-17:30", "17:30-18:30")

foreach ($sheetName in @("Section_A", "Section_B")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $startRow = 8
    for ($i = 0; $i -lt $newTimes.Length; $i++) {
        $row = $startRow + $i

        $ws.Cells.Item($row, 1).Value = $newTimes[$i]
        $ws.Cells.Item($row, 2).Value = "Free"
        $ws.Cells.Item($row, 3).Value = "Free"
        $ws.Cells.Item($row, 4).Value = "Free"
        $ws.Cells.Item($row, 5).Value = "Free"
        $ws.Cells.Item($row, 6).Value = "Free"

        # Copy the time-column style (bordered / bold / centered) from the
        # existing A2 cell onto the new A column cell.
        $ws.Range("A2").Copy()
        $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    }
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Course_Summary sheet: split the single "Credits"/"Instructor" columns into
# Lectures/Week, Tutorials/Week, Total Credits, Instructor - since tutorials
# are now a separate weekly component that is tracked on its own.
# ---------------------------------------------------------------------------
$cs = $wb.Worksheets.Item("Course_Summary")

$cs.Range("E1").Value = "Lectures/Week"
$cs.Range("F1").Value = "Tutorials/Week"
$cs.Range("G1").Value = "Total Credits"
$cs.Range("H1").Value = "Instructor"

# Give the two new header cells (G1, H1) the same header styling as F1.
$cs.Range("F1").Copy()
$cs.Range("G1").PasteSpecial(-4122)
$cs.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$instructors = @{ 2 = "Dr. Raghav Menon"; 3 = "Dr. Sneha Rao"; 4 = "Dr. Neel Patel" }

foreach ($row in 2..4) {
    $cs.Cells.Item($row, 6).Value = 0               # F: Tutorials/Week
    $cs.Cells.Item($row, 7).Value = 4                # G: Total Credits
    $cs.Cells.Item($row, 8).Value = $instructors[$row]  # H: Instructor
}
